$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy the row style (fill/border/font/alignment) from column C to column D
#     for every data row, before writing values (PasteSpecial formats only). ---
$ws.Range("C2").Copy()
$ws.Range("D2").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("D3").PasteSpecial(-4122)
$ws.Range("C4").Copy()
$ws.Range("D4").PasteSpecial(-4122)
$ws.Range("C5").Copy()
$ws.Range("D5").PasteSpecial(-4122)
$ws.Range("C6").Copy()
$ws.Range("D6").PasteSpecial(-4122)
$ws.Range("C7").Copy()
$ws.Range("D7").PasteSpecial(-4122)
$ws.Range("C8").Copy()
$ws.Range("D8").PasteSpecial(-4122)
$ws.Range("C9").Copy()
$ws.Range("D9").PasteSpecial(-4122)
$ws.Range("C10").Copy()
$ws.Range("D10").PasteSpecial(-4122)
$ws.Range("C11").Copy()
$ws.Range("D11").PasteSpecial(-4122)
$ws.Range("C12").Copy()
$ws.Range("D12").PasteSpecial(-4122)
$ws.Range("C13").Copy()
$ws.Range("D13").PasteSpecial(-4122)

# --- Header: new column D "canonical SMILES" ---
$ws.Range("D2").Value = "canonical SMILES"

# --- New, simplified (non-stereo) canonical SMILES for each microstate. ---
# Rows where a brand-new SMILES string is introduced in column D (the
# stereo-bond markers '/' and '\' are dropped relative to column C).
$ws.Range("D3").Value  = "COc1cccc(c1)[NH+]=c2c3ccccc3[nH]cn2"
$ws.Range("D4").Value  = "COc1cccc(c1)N=c2c3ccccc3[nH]cn2"
$ws.Range("D10").Value = "COc1cccc(c1)N=c2c3ccccc3nc[nH]2"

# Rows where column C was already a plain (non-stereo) SMILES, so column D
# simply repeats the same value as column C.
$ws.Range("D5").Value  = $ws.Range("C5").Value()
$ws.Range("D6").Value  = $ws.Range("C6").Value()
$ws.Range("D7").Value  = $ws.Range("C7").Value()
$ws.Range("D8").Value  = $ws.Range("C8").Value()
$ws.Range("D9").Value  = $ws.Range("C9").Value()
$ws.Range("D11").Value = $ws.Range("C11").Value()
$ws.Range("D12").Value = $ws.Range("C12").Value()
$ws.Range("D13").Value = $ws.Range("C13").Value()

# --- New column D width (target OOXML width 36.85546875 chars; the closest
#     value this engine's pixel-quantized ColumnWidth model can reach). ---
$ws.Columns.Item(4).ColumnWidth = 36.0

Write-Host "done"
